$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 58 (existing rows 58-76 shift down to 61-79)
$ws.Range("A58:A60").EntireRow.Insert()

# Common values shared by these "Chirimoya" records
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "$/kilo (en caja de 15 kilos)"
$origen      = "Provincia de Limarí"
$fecha       = 44524

# New row 58: Especial
$r = 58
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $tipo
$ws.Cells.Item($r,7).Value = $productoId
$ws.Cells.Item($r,8).Value = $producto
$ws.Cells.Item($r,9).Value = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 300
$ws.Cells.Item($r,14).Value = 1700
$ws.Cells.Item($r,15).Value = 1800
$ws.Cells.Item($r,16).Value = 1750
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 1750
$ws.Cells.Item($r,20).Value = 1

# New row 59: Primera
$r = 59
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $tipo
$ws.Cells.Item($r,7).Value = $productoId
$ws.Cells.Item($r,8).Value = $producto
$ws.Cells.Item($r,9).Value = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 300
$ws.Cells.Item($r,14).Value = 1400
$ws.Cells.Item($r,15).Value = 1500
$ws.Cells.Item($r,16).Value = 1450
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 1450
$ws.Cells.Item($r,20).Value = 1

# New row 60: Segunda
$r = 60
$ws.Cells.Item($r,1).Value = $mercadoId
$ws.Cells.Item($r,2).Value = $mercado
$ws.Cells.Item($r,3).Value = $region
$ws.Cells.Item($r,4).Value = $fecha
$ws.Cells.Item($r,5).Value = $codreg
$ws.Cells.Item($r,6).Value = $tipo
$ws.Cells.Item($r,7).Value = $productoId
$ws.Cells.Item($r,8).Value = $producto
$ws.Cells.Item($r,9).Value = $categoriaId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = $variedad
$ws.Cells.Item($r,12).Value = "Segunda"
$ws.Cells.Item($r,13).Value = 240
$ws.Cells.Item($r,14).Value = 1100
$ws.Cells.Item($r,15).Value = 1200
$ws.Cells.Item($r,16).Value = 1150
$ws.Cells.Item($r,17).Value = $unidad
$ws.Cells.Item($r,18).Value = $origen
$ws.Cells.Item($r,19).Value = 1150
$ws.Cells.Item($r,20).Value = 1
